# edit.ps1
# Applies the 2026-02-16 Betfair odds update:
#   - updates many Back/Lay odds values throughout the sheet
#   - swaps the two "Dutch Eerste Divisie" fixtures (rows 15/16) together with their odds
#   - inserts a new "Colombian Primera B" fixture as row 24 (pushing the following three
#     rows down by one)
#   - appends a new "Uruguayan Primera Division" fixture as the final row (28)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Step 1: insert a new row at position 24 for the Colombian Primera B
# fixture. This shifts the old rows 24-26 down to rows 25-27.
# ------------------------------------------------------------------
$ws.Rows.Item(24).Insert()

# Date/Time columns store plain text (e.g. "2026-02-16"), but assigning
# such a string straight to .Value makes Excel auto-convert it to a real
# date/time serial number. Copy/PasteSpecial(values) from a neighbouring
# cell that already holds the same text keeps it as plain text.
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B24").PasteSpecial(-4163) | Out-Null
$ws.Range("C25").Copy() | Out-Null
$ws.Range("C24").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Step 2: populate the rest of the newly inserted row 24
# (League, Home, Away, and all the odds columns)
# ------------------------------------------------------------------
$ws.Cells.Item(24, 1).Value = "Colombian Primera B"
$ws.Cells.Item(24, 4).Value = "Tigres FC Zipaquira"
$ws.Cells.Item(24, 5).Value = "Independiente Yumbo"
$ws.Cells.Item(24, 6).Value = 1.04
$ws.Cells.Item(24, 7).Value = 1000
$ws.Cells.Item(24, 8).Value = 1.04
$ws.Cells.Item(24, 9).Value = 1000
$ws.Cells.Item(24, 10).Value = 1.02
$ws.Cells.Item(24, 11).Value = 1000
$ws.Cells.Item(24, 12).Value = 1.01
$ws.Cells.Item(24, 13).Value = 1.01
$ws.Cells.Item(24, 14).Value = 1.25
$ws.Cells.Item(24, 15).Value = 1.02
$ws.Cells.Item(24, 16).Value = 1.24
$ws.Cells.Item(24, 17).Value = 1.43
$ws.Cells.Item(24, 18).Value = 1.13
$ws.Cells.Item(24, 19).Value = 1.43
$ws.Cells.Item(24, 20).Value = 1.01
$ws.Cells.Item(24, 21).Value = 1.01
$ws.Cells.Item(24, 22).Value = 1.01
$ws.Cells.Item(24, 23).Value = 1.01
$ws.Cells.Item(24, 24).Value = 1000
$ws.Cells.Item(24, 25).Value = 1000
$ws.Cells.Item(24, 26).Value = 1000
$ws.Cells.Item(24, 27).Value = 1000
$ws.Cells.Item(24, 28).Value = 1000
$ws.Cells.Item(24, 29).Value = 1000
$ws.Cells.Item(24, 30).Value = 1000
$ws.Cells.Item(24, 31).Value = 1000
$ws.Cells.Item(24, 32).Value = 1000
$ws.Cells.Item(24, 33).Value = 1000
$ws.Cells.Item(24, 34).Value = 1000
$ws.Cells.Item(24, 35).Value = 1000
$ws.Cells.Item(24, 36).Value = 1000
$ws.Cells.Item(24, 37).Value = 1000
$ws.Cells.Item(24, 38).Value = 1000
$ws.Cells.Item(24, 39).Value = 1000
$ws.Cells.Item(24, 40).Value = 1000
$ws.Cells.Item(24, 41).Value = 1000

# ------------------------------------------------------------------
# Step 3: apply the updated odds / team-name values for the rows that
# keep their original position (1-14, 17-23), the two swapped Dutch
# Eerste Divisie fixtures (15-16), and the three fixtures that were
# shifted down by the row insertion above (now rows 25-27).
# ------------------------------------------------------------------
$ws.Cells.Item(2, 10).Value = 1.09
$ws.Cells.Item(2, 14).Value = 1.02
$ws.Cells.Item(2, 15).Value = 1.01
$ws.Cells.Item(2, 17).Value = 1.01
$ws.Cells.Item(2, 19).Value = 1.01
$ws.Cells.Item(3, 6).Value = 5.1
$ws.Cells.Item(3, 7).Value = 6.4
$ws.Cells.Item(3, 8).Value = 1.82
$ws.Cells.Item(3, 9).Value = 1.84
$ws.Cells.Item(3, 11).Value = 3.8
$ws.Cells.Item(3, 16).Value = 1.66
$ws.Cells.Item(3, 17).Value = 2.26
$ws.Cells.Item(4, 6).Value = 1.74
$ws.Cells.Item(4, 8).Value = 4.3
$ws.Cells.Item(4, 9).Value = 7.4
$ws.Cells.Item(4, 10).Value = 2.82
$ws.Cells.Item(4, 11).Value = 5.2
$ws.Cells.Item(4, 16).Value = 1.45
$ws.Cells.Item(4, 17).Value = 2.24
$ws.Cells.Item(5, 6).Value = 1.4
$ws.Cells.Item(7, 6).Value = 1.24
$ws.Cells.Item(7, 7).Value = 2.38
$ws.Cells.Item(7, 8).Value = 3.5
$ws.Cells.Item(7, 10).Value = 3.4
$ws.Cells.Item(7, 11).Value = 980
$ws.Cells.Item(8, 16).Value = 1.81
$ws.Cells.Item(9, 8).Value = 3.65
$ws.Cells.Item(9, 17).Value = 1.69
$ws.Cells.Item(10, 6).Value = 3.1
$ws.Cells.Item(10, 7).Value = 4.4
$ws.Cells.Item(10, 8).Value = 2.04
$ws.Cells.Item(10, 10).Value = 3.6
$ws.Cells.Item(10, 11).Value = 5.7
$ws.Cells.Item(11, 6).Value = 2.38
$ws.Cells.Item(11, 8).Value = 2.78
$ws.Cells.Item(11, 10).Value = 2.98
$ws.Cells.Item(11, 11).Value = 4.6
$ws.Cells.Item(11, 16).Value = 1.76
$ws.Cells.Item(11, 17).Value = 1.78
$ws.Cells.Item(13, 10).Value = 5
$ws.Cells.Item(14, 7).Value = 3.35
$ws.Cells.Item(14, 8).Value = 2.3
$ws.Cells.Item(14, 10).Value = 3.75
$ws.Cells.Item(14, 11).Value = 6.8
$ws.Cells.Item(14, 16).Value = 2.38
$ws.Cells.Item(21, 6).Value = 2.9
$ws.Cells.Item(21, 7).Value = 3.25
$ws.Cells.Item(21, 9).Value = 2.74
$ws.Cells.Item(21, 16).Value = 1.87
$ws.Cells.Item(21, 17).Value = 1.95
$ws.Cells.Item(22, 7).Value = 3.25
$ws.Cells.Item(22, 9).Value = 3.8
$ws.Cells.Item(22, 10).Value = 2.64
$ws.Cells.Item(22, 16).Value = 1.58
$ws.Cells.Item(22, 17).Value = 2.24
$ws.Cells.Item(23, 8).Value = 4.1
$ws.Cells.Item(23, 15).Value = 1.55
$ws.Cells.Item(23, 16).Value = 1.54
$ws.Cells.Item(23, 17).Value = 2.74
$ws.Cells.Item(23, 26).Value = 28
$ws.Cells.Item(23, 31).Value = 1000
$ws.Cells.Item(23, 40).Value = 32
$ws.Cells.Item(15, 4).Value = "Jong FC Utrecht"
$ws.Cells.Item(15, 5).Value = "Willem II"
$ws.Cells.Item(15, 6).Value = 2.32
$ws.Cells.Item(15, 7).Value = 3.2
$ws.Cells.Item(15, 8).Value = 2.16
$ws.Cells.Item(15, 9).Value = 2.9
$ws.Cells.Item(15, 10).Value = 3.8
$ws.Cells.Item(15, 11).Value = 10.5
$ws.Cells.Item(15, 16).Value = 2.8
$ws.Cells.Item(15, 17).Value = 1.4
$ws.Cells.Item(16, 4).Value = "Jong Ajax Amsterdam"
$ws.Cells.Item(16, 5).Value = "Cambuur Leeuwarden"
$ws.Cells.Item(16, 6).Value = 5.7
$ws.Cells.Item(16, 7).Value = 7.2
$ws.Cells.Item(16, 8).Value = 1.52
$ws.Cells.Item(16, 9).Value = 1.62
$ws.Cells.Item(16, 10).Value = 4.9
$ws.Cells.Item(16, 11).Value = 6
$ws.Cells.Item(16, 16).Value = 3
$ws.Cells.Item(25, 17).Value = 1.36
$ws.Cells.Item(25, 18).Value = 2.04
$ws.Cells.Item(25, 19).Value = 1.92
$ws.Cells.Item(25, 21).Value = 2.64
$ws.Cells.Item(25, 24).Value = 48
$ws.Cells.Item(25, 28).Value = 46
$ws.Cells.Item(25, 39).Value = 1000
$ws.Cells.Item(25, 40).Value = 460
$ws.Cells.Item(26, 7).Value = 3.6
$ws.Cells.Item(26, 8).Value = 2.48
$ws.Cells.Item(27, 6).Value = 2.32
$ws.Cells.Item(27, 9).Value = 4.6
$ws.Cells.Item(27, 10).Value = 2.7
$ws.Cells.Item(27, 11).Value = 3.05

# ------------------------------------------------------------------
# Step 4: append the new row 28 for the Uruguayan Primera Division
# fixture.
# ------------------------------------------------------------------
$ws.Range("B27").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4163) | Out-Null
$ws.Range("C27").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

$ws.Cells.Item(28, 1).Value = "Uruguayan Primera Division"
$ws.Cells.Item(28, 4).Value = "Albion FC"
$ws.Cells.Item(28, 5).Value = "Cerro"
$ws.Cells.Item(28, 6).Value = 2.16
$ws.Cells.Item(28, 7).Value = 2.54
$ws.Cells.Item(28, 8).Value = 3.8
$ws.Cells.Item(28, 9).Value = 4.9
$ws.Cells.Item(28, 10).Value = 2.92
$ws.Cells.Item(28, 11).Value = 3.35
$ws.Cells.Item(28, 12).Value = 0
$ws.Cells.Item(28, 13).Value = 0
$ws.Cells.Item(28, 14).Value = 0
$ws.Cells.Item(28, 15).Value = 0
$ws.Cells.Item(28, 16).Value = 1.51
$ws.Cells.Item(28, 17).Value = 2.38
$ws.Cells.Item(28, 18).Value = 0
$ws.Cells.Item(28, 19).Value = 0
$ws.Cells.Item(28, 20).Value = 0
$ws.Cells.Item(28, 21).Value = 0
$ws.Cells.Item(28, 22).Value = 0
$ws.Cells.Item(28, 23).Value = 0
$ws.Cells.Item(28, 24).Value = 0
$ws.Cells.Item(28, 25).Value = 0
$ws.Cells.Item(28, 26).Value = 0
$ws.Cells.Item(28, 27).Value = 0
$ws.Cells.Item(28, 28).Value = 0
$ws.Cells.Item(28, 29).Value = 0
$ws.Cells.Item(28, 30).Value = 0
$ws.Cells.Item(28, 31).Value = 0
$ws.Cells.Item(28, 32).Value = 0
$ws.Cells.Item(28, 33).Value = 0
$ws.Cells.Item(28, 34).Value = 0
$ws.Cells.Item(28, 35).Value = 0
$ws.Cells.Item(28, 36).Value = 0
$ws.Cells.Item(28, 37).Value = 0
$ws.Cells.Item(28, 38).Value = 0
$ws.Cells.Item(28, 39).Value = 0
$ws.Cells.Item(28, 40).Value = 0
$ws.Cells.Item(28, 41).Value = 0

Write-Output "Edit applied successfully."
